$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter ..." footer line.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter*Salvar em pdf*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # The paragraph immediately before it is the blank separator paragraph,
    # and the paragraph immediately after it is the "© 2020 ..." footer line.
    # Removing all three (separator + "Ver no Jupiter..." + "© 2020...") collapses
    # back onto the blank paragraph that precedes the trailing page break.
    $startPara = $d.Paragraphs.Item($target - 1)
    $endPara = $d.Paragraphs.Item($target + 1)

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
